# Regenerate the handoff status report:
#  - the first "source" file (old UUID) got a new generated UUID + transform hash,
#    and its handoff timestamps moved forward a minute
#  - the failed-transform row for the old second source file is gone entirely
#    (that source file no longer produced a failed transform), so that row is
#    removed from every sheet and the ".localization-config" row moves up.

$oldUuid = "63f6109c-1561-4875-8e6c-eb30ed2f8749"
$newUuid = "b4327c94-b673-4614-b729-626121798d22"
$oldHash = "b8fdec19845ad76d90b507ad86db1ab590d830d1"
$newHash = "3fa6b060e40ccf0b67a8f0cf860aa2e12c69d407"
$oldDtZh = "2016-01-22 02:51:03"
$newDtZh = "2016-01-22 02:52:03"
$oldDtDe = "2016-01-22 02:51:17"
$newDtDe = "2016-01-22 02:52:17"

$removedMdName = "4dfd6a95-976c-4246-a326-27ccde1a8613.md"

$wb = $excel.ActiveWorkbook

# ---------- Sheet "Overview" ----------
$ws = $wb.Worksheets.Item("Overview")

# Drop every hyperlink on the sheet; we'll recreate the ones that remain below.
$ws.Range("A1").Hyperlinks.Delete()

# Remove the row describing the now-gone "4dfd6a95...md" source (row 3); the
# ".localization-config" row that used to be row 4 slides up to row 3.
$ws.Rows.Item(3).Delete()

# Update the remaining source file's name wherever it appears.
$ws.Range("A2").Value = $newUuid + ".md"

# Re-create the two remaining hyperlinks.
$h1 = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/e2e/" + $newUuid + ".md", "", "", $newUuid + ".md")
$h2 = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/.localization-config", "", "", ".localization-config")

# ---------- Sheet "zh-cn" ----------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newUuid + ".md"
$ws.Range("C2").Value = $newUuid + "." + $newHash + ".zh-cn.xlf"
$ws.Range("D2").Value = $newDtZh

$h1 = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/e2e/" + $newUuid + ".md", "", "", $newUuid + ".md")
$h2 = $ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/616289dbce928d04d4df286e4025d567c8a523ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/" + $newUuid + "." + $newHash + ".zh-cn.xlf", "", "", $newUuid + "." + $newHash + ".zh-cn.xlf")
$h3 = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/.localization-config", "", "", ".localization-config")

# ---------- Sheet "de-de" ----------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newUuid + ".md"
$ws.Range("C2").Value = $newUuid + "." + $newHash + ".de-de.xlf"
$ws.Range("D2").Value = $newDtDe

$h1 = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/e2e/" + $newUuid + ".md", "", "", $newUuid + ".md")
$h2 = $ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1cce06d6afb76055398998b5b9861761c15c1930/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/" + $newUuid + "." + $newHash + ".de-de.xlf", "", "", $newUuid + "." + $newHash + ".de-de.xlf")
$h3 = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f67c27c606033cbbc3759d370e27890c63134c91/.localization-config", "", "", ".localization-config")

$wb.Save()
